$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "Ccl4"
$ws.Range("C2").Value = "Ackr2"
$ws.Range("G2").Value = 0.8774723333333333
$ws.Range("H2").Value = 2.6324169999999998
$ws.Range("I2").Value = 0.00070883934342592708
$ws.Range("J2").Value = 0.00070883934342592708
$ws.Range("Q2").Value = 0.22991530078
$ws.Range("R2").Value = 2.0692377070200001
$ws.Range("S2").Value = 0.00070883934342592708
$ws.Range("T2").Value = 0.00070883934342592708

# Row 3
$ws.Range("B3").Value = "Ccl4"
$ws.Range("C3").Value = "Ackr2"
$ws.Range("I3").Value = 0.00006864123866041489
$ws.Range("J3").Value = 0.00006864123866041489
$ws.Range("S3").Value = 0.00006864123866041489
$ws.Range("T3").Value = 0.00006864123866041489

# Row 4
$ws.Range("B4").Value = "Ccl4"
$ws.Range("C4").Value = "Ackr2"
$ws.Range("G4").Value = 71.071772333333342
$ws.Range("H4").Value = 213.215317
$ws.Range("I4").Value = 0.057413170219851553
$ws.Range("J4").Value = 0.057413170219851539
$ws.Range("Q4").Value = 18.62222578678
$ws.Range("R4").Value = 167.60003208102
$ws.Range("S4").Value = 0.057413170219851553
$ws.Range("T4").Value = 0.057413170219851539

# Row 5 (Sending cluster becomes MuSCs)
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Ccl4"
$ws.Range("C5").Value = "Ackr2"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.33333333333333331
$ws.Range("G5").Value = 0.035969333333333332
$ws.Range("H5").Value = 0.107908
$ws.Range("I5").Value = 0.00002905673222380989
$ws.Range("J5").Value = 0.00002905673222380989
$ws.Range("Q5").Value = 0.009424684719999998
$ws.Range("R5").Value = 0.084822162480000005
$ws.Range("S5").Value = 0.00002905673222380989
$ws.Range("T5").Value = 0.00002905673222380989

# Row 6 (Sending cluster becomes Neutrophils)
$ws.Range("A6").Value = "Neutrophils"
$ws.Range("B6").Value = "Ccl4"
$ws.Range("C6").Value = "Ackr2"
$ws.Range("G6").Value = 1112.7588499999999
$ws.Range("H6").Value = 3338.27655
$ws.Range("I6").Value = 0.8989084016233635
$ws.Range("J6").Value = 0.89890840162336338
$ws.Range("Q6").Value = 291.56507387699997
$ws.Range("R6").Value = 2624.0856648929998
$ws.Range("S6").Value = 0.8989084016233635
$ws.Range("T6").Value = 0.89890840162336338

# Row 7 (new row, Resolving-Mac)
$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("B7").Value = "Ccl4"
$ws.Range("C7").Value = "Ackr2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 53.071120333333333
$ws.Range("H7").Value = 159.21336099999999
$ws.Range("I7").Value = 0.04287189084247485
$ws.Range("J7").Value = 0.042871890842474843
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.66666666666666663
$ws.Range("M7").Value = 0.26201999999999998
$ws.Range("N7").Value = 0.78605999999999998
$ws.Range("O7").Value = 1
$ws.Range("P7").Value = 1
$ws.Range("Q7").Value = 13.905694949740001
$ws.Range("R7").Value = 125.15125454766
$ws.Range("S7").Value = 0.04287189084247485
$ws.Range("T7").Value = 0.042871890842474843
